$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 883.2632
$ws.Range("I19").Value = 330.16666
$ws.Range("J19").Value = 1138.5385
$ws.Range("K19").Value = 330.16666
$ws.Range("L19").Value = 1138.5385
$ws.Range("M19").Value = -155.16666
$ws.Range("N19").Value = -1488.5385
$ws.Range("H37").Value = 248.4
$ws.Range("I37").Value = 248.4
$ws.Range("K37").Value = 745.2
$ws.Range("M37").Value = -619.2
$ws.Range("H115").Value = 185
$ws.Range("I115").Value = 185
$ws.Range("K115").Value = 555
$ws.Range("M115").Value = 1012
$ws.Range("H135").Value = 1283.0667
$ws.Range("I135").Value = 883.1667
$ws.Range("J135").Value = 2882.6667
$ws.Range("K135").Value = 7948.5003
$ws.Range("L135").Value = 25944.0003
$ws.Range("M135").Value = -5413.5003
$ws.Range("N135").Value = -31014.0003
$ws.Range("H137").Value = 11849.1
$ws.Range("I137").Value = 1998.8572
$ws.Range("J137").Value = 34833
$ws.Range("K137").Value = 5996.571599999999
$ws.Range("L137").Value = 104499
$ws.Range("M137").Value = -3446.571599999999
$ws.Range("N137").Value = -109599
$ws.Range("H138").Value = 29443.975
$ws.Range("I138").Value = 61647.47
$ws.Range("J138").Value = 5641.391
$ws.Range("K138").Value = 184942.41
$ws.Range("L138").Value = 16924.173
$ws.Range("M138").Value = -179802.41
$ws.Range("N138").Value = -27204.173

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4392.7104
$ws.Range("I32").Value = 4392.7104
$ws.Range("K32").Value = 4392.7104
$ws.Range("M32").Value = -4105.7104
$ws.Range("H45").Value = 123788.11
$ws.Range("I45").Value = 123788.11
$ws.Range("K45").Value = 123788.11
$ws.Range("M45").Value = -123411.11
$ws.Range("H74").Value = 11696.177
$ws.Range("I74").Value = 1402.9166
$ws.Range("K74").Value = 1402.9166
$ws.Range("M74").Value = -528.9166
$ws.Range("H77").Value = 11696.177
$ws.Range("I77").Value = 1402.9166
$ws.Range("K77").Value = 7014.583000000001
$ws.Range("M77").Value = -2646.583000000001
$ws.Range("H122").Value = 1910.069
$ws.Range("I122").Value = 1950.96
$ws.Range("K122").Value = 5852.88
$ws.Range("M122").Value = -3402.88

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 13222
$ws.Range("J106").Value = 13222
$ws.Range("L106").Value = 13222
$ws.Range("N106").Value = -15746
$ws.Range("H107").Value = 2190.1667
$ws.Range("I107").Value = 1878.4783
$ws.Range("J107").Value = 3214.2856
$ws.Range("K107").Value = 1878.4783
$ws.Range("L107").Value = 3214.2856
$ws.Range("M107").Value = 41.52170000000001
$ws.Range("N107").Value = -7054.2856

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2465.25
$ws.Range("I16").Value = 1524.8
$ws.Range("J16").Value = 4032.6667
$ws.Range("K16").Value = 1524.8
$ws.Range("L16").Value = 4032.6667
$ws.Range("M16").Value = -1237.8
$ws.Range("N16").Value = -4606.6667
$ws.Range("H31").Value = 23340.084
$ws.Range("I31").Value = 32384.697
$ws.Range("J31").Value = 3441.9333
$ws.Range("K31").Value = 32384.697
$ws.Range("L31").Value = 3441.9333
$ws.Range("M31").Value = -32089.697
$ws.Range("N31").Value = -4031.9333
$ws.Range("H34").Value = 23340.084
$ws.Range("I34").Value = 32384.697
$ws.Range("J34").Value = 3441.9333
$ws.Range("K34").Value = 32384.697
$ws.Range("L34").Value = 3441.9333
$ws.Range("M34").Value = -32182.697
$ws.Range("N34").Value = -3845.9333
$ws.Range("H58").Value = 3264.6938
$ws.Range("I58").Value = 3052.0857
$ws.Range("J58").Value = 3796.2144
$ws.Range("K58").Value = 3052.0857
$ws.Range("L58").Value = 3796.2144
$ws.Range("M58").Value = -2849.0857
$ws.Range("N58").Value = -4202.2144
$ws.Range("H99").Value = 3575.5454
$ws.Range("I99").Value = 4111.3335
$ws.Range("J99").Value = 3374.625
$ws.Range("K99").Value = 4111.3335
$ws.Range("L99").Value = 3374.625
$ws.Range("M99").Value = -2613.3335
$ws.Range("N99").Value = -6370.625
$ws.Range("H107").Value = 970.41174
$ws.Range("I107").Value = 1087.2858
$ws.Range("J107").Value = 888.6
$ws.Range("K107").Value = 1087.2858
$ws.Range("L107").Value = 888.6
$ws.Range("M107").Value = 832.7141999999999
$ws.Range("N107").Value = -4728.6
$ws.Range("H113").Value = 2465.25
$ws.Range("I113").Value = 1524.8
$ws.Range("J113").Value = 4032.6667
$ws.Range("K113").Value = 1524.8
$ws.Range("L113").Value = 4032.6667
$ws.Range("M113").Value = 645.2
$ws.Range("N113").Value = -8372.6667
$ws.Range("H114").Value = 22500
$ws.Range("J114").Value = 22500
$ws.Range("L114").Value = 22500
$ws.Range("N114").Value = -31178
$ws.Range("H126").Value = 3575.5454
$ws.Range("I126").Value = 4111.3335
$ws.Range("J126").Value = 3374.625
$ws.Range("K126").Value = 12334.0005
$ws.Range("L126").Value = 10123.875
$ws.Range("M126").Value = -9864.000499999998
$ws.Range("N126").Value = -15063.875
$ws.Range("H132").Value = 172927.56
$ws.Range("I132").Value = 241833.58
$ws.Range("J132").Value = 12146.889
$ws.Range("K132").Value = 725500.74
$ws.Range("L132").Value = 36440.667
$ws.Range("M132").Value = -722970.74
$ws.Range("N132").Value = -41500.667
$ws.Range("H136").Value = 3264.6938
$ws.Range("I136").Value = 3052.0857
$ws.Range("J136").Value = 3796.2144
$ws.Range("K136").Value = 9156.257100000001
$ws.Range("L136").Value = 11388.6432
$ws.Range("M136").Value = -6606.257100000001
$ws.Range("N136").Value = -16488.6432

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 3
$ws.Range("I47").Value = 3
$ws.Range("K47").Value = 9
$ws.Range("M47").Value = 422
$ws.Range("H87").Value = 15571.611
$ws.Range("I87").Value = 10023.385
$ws.Range("K87").Value = 30070.155
$ws.Range("M87").Value = -28822.155
$ws.Range("H90").Value = 15571.611
$ws.Range("I90").Value = 10023.385
$ws.Range("K90").Value = 90210.465
$ws.Range("M90").Value = -83970.465
$ws.Range("H113").Value = 1004.38464
$ws.Range("I113").Value = 400
$ws.Range("J113").Value = 1054.75
$ws.Range("K113").Value = 1200
$ws.Range("L113").Value = 3164.25
$ws.Range("M113").Value = 970
$ws.Range("N113").Value = -7504.25
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").ClearContents()
$ws.Range("H122").Value = 1786
$ws.Range("J122").Value = 1136.875
$ws.Range("L122").Value = 10231.875
$ws.Range("N122").Value = -15131.875
$ws.Range("H129").Value = 834.2
$ws.Range("I129").Value = 542.6667
$ws.Range("J129").Value = 2000.3334
$ws.Range("K129").Value = 1628.0001
$ws.Range("L129").Value = 6001.0002
$ws.Range("M129").Value = 3371.9999
$ws.Range("N129").Value = -16001.0002

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 57177
$ws.Range("J100").Value = 57177
$ws.Range("L100").Value = 57177
$ws.Range("N100").Value = -59341
$ws.Range("H122").Value = 3541.4285
$ws.Range("I122").Value = 2932
$ws.Range("J122").Value = 3998.5
$ws.Range("K122").Value = 8796
$ws.Range("L122").Value = 11995.5
$ws.Range("M122").Value = -6346
$ws.Range("N122").Value = -16895.5
$ws.Range("H132").Value = 1745.4517
$ws.Range("I132").Value = 1570.3
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 4710.9
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -2180.9
$ws.Range("N132").Value = -26060

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2821.4443
$ws.Range("I93").Value = 3024.25
$ws.Range("J93").Value = 1199
$ws.Range("K93").Value = 3024.25
$ws.Range("L93").Value = 1199
$ws.Range("M93").Value = -1776.25
$ws.Range("N93").Value = -3695
$ws.Range("H107").Value = 6729.5
$ws.Range("I107").Value = 6729.5
$ws.Range("K107").Value = 6729.5
$ws.Range("M107").Value = -4809.5
$ws.Range("H132").Value = 2859.375
$ws.Range("I132").Value = 2699.0732
$ws.Range("J132").Value = 3798.2856
$ws.Range("K132").Value = 8097.219599999999
$ws.Range("L132").Value = 11394.8568
$ws.Range("M132").Value = -5567.219599999999
$ws.Range("N132").Value = -16454.8568
$ws.Range("H136").Value = 62788.234
$ws.Range("I136").Value = 80877.38
$ws.Range("J136").Value = 3998.5
$ws.Range("K136").Value = 242632.14
$ws.Range("L136").Value = 11995.5
$ws.Range("M136").Value = -240082.14
$ws.Range("N136").Value = -17095.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 17151.334
$ws.Range("J45").Value = 17151.334
$ws.Range("L45").Value = 17151.334
$ws.Range("N45").Value = -18133.334
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 2063.5312
$ws.Range("I132").Value = 1900.7307
$ws.Range("J132").Value = 2769
$ws.Range("K132").Value = 5702.1921
$ws.Range("L132").Value = 8307
$ws.Range("M132").Value = -3172.1921
$ws.Range("N132").Value = -13367
$ws.Range("H136").Value = 2510
$ws.Range("I136").Value = 2495.4
$ws.Range("J136").Value = 2537.375
$ws.Range("K136").Value = 7486.200000000001
$ws.Range("L136").Value = 7612.125
$ws.Range("M136").Value = -4936.200000000001
$ws.Range("N136").Value = -12712.125
